$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.439.76"
$ws.Range("E2").Value = "  +4.67%  "
$ws.Range("D3").Value = "2.752.47"
$ws.Range("E3").Value = "  +4.66%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "115.91"
$ws.Range("E5").Value = "  +3.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "332.41"
$ws.Range("E6").Value = "  +2.99%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.539"
$ws.Range("E7").Value = "  +2.47%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.574"
$ws.Range("E9").Value = "  +5.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.60"
$ws.Range("E10").Value = "  +4.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0856"
$ws.Range("E11").Value = "  +5.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.22"
$ws.Range("E12").Value = "  +2.47%  "
$ws.Range("E13").Value = "  +2.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.63"
$ws.Range("E14").Value = "  +5.20%  "
$ws.Range("D15").Value = "3.175.03"
$ws.Range("E15").Value = "  +4.56%  "
$ws.Range("D16").Value = "2.773.51"
$ws.Range("E16").Value = "  +5.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.886"
$ws.Range("E17").Value = "  +3.19%  "
$ws.Range("D18").Value = "51.369.01"
$ws.Range("E18").Value = "  +4.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.20"
$ws.Range("E19").Value = "  +5.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.50"
$ws.Range("E20").Value = "  +4.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.85"
$ws.Range("E21").Value = "  +2.26%  "
$ws.Range("D22").Value = "0.0₃0975"
$ws.Range("E22").Value = "  +3.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "278.38"
$ws.Range("E23").Value = "  +3.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.62"
$ws.Range("E24").Value = "  +1.59%  "
$ws.Range("E25").Value = "  +4.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.79"
$ws.Range("E26").Value = "  +2.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.18"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("E29").Value = "  -0.58%  "
$ws.Range("E30").Value = "  +1.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.25"
$ws.Range("E31").Value = "  +0.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.88"
$ws.Range("E32").Value = "  +0.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.55"
$ws.Range("E33").Value = "  +1.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0825"
$ws.Range("E34").Value = "  +3.01%  "
$ws.Range("B35").Value = "Celestia"
$ws.Range("C35").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.10"
$ws.Range("E35").Value = "  +0.38%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.09"
$ws.Range("E37").Value = "  +2.50%  "
$ws.Range("E38").Value = "  +0.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.25"
$ws.Range("E39").Value = "  +4.03%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0351"
$ws.Range("E40").Value = "  +10.57%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "126.73"
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "23.15"
$ws.Range("E42").Value = "  +4.29%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.114"
$ws.Range("E43").Value = "  +2.86%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.29"
$ws.Range("E44").Value = "  +7.81%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.44"
$ws.Range("E45").Value = "  +13.32%  "
$ws.Range("D46").Value = "2.089.73"
$ws.Range("E46").Value = "  +1.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.33"
$ws.Range("E47").Value = "  +3.93%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.52"
$ws.Range("E49").Value = "  +6.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.98"
$ws.Range("E50").Value = "  +0.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "59.94"
$ws.Range("E51").Value = "  +2.08%  "
